$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the remaining data row (row 2) with the new date serial and value
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 11.51866285751828

# Remove the now-obsolete trailing rows (3 through 17) entirely so the
# sheet's used range / dimension shrinks back down to A1:B2
$ws.Range("A3:B17").EntireRow.Delete()
